$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 previously only had A36 populated (with style) and no B36 cell at all.
# Copy the date-format from the row above so the new B36 cell matches the
# existing B-column date formatting (style index carried over, not a new one).
$ws.Range("B35").Copy()
$ws.Range("B36").PasteSpecial(-4122)

# The whole shift schedule (weekday label in column A, date serial in column B)
# moved back by one weekday / 31 days, and a new trailing row (36) was added to
# complete the final week.
$ws.Range("A6").Value = "Donnerstag"
$ws.Range("B6").Value2 = 42216
$ws.Range("A7").Value = "Freitag"
$ws.Range("B7").Value2 = 42217
$ws.Range("A8").Value = "Samstag"
$ws.Range("B8").Value2 = 42218
$ws.Range("A9").Value = "Sonntag"
$ws.Range("B9").Value2 = 42219
$ws.Range("A10").Value = "Montag"
$ws.Range("B10").Value2 = 42220
$ws.Range("A11").Value = "Dienstag"
$ws.Range("B11").Value2 = 42221
$ws.Range("A12").Value = "Mittwoch"
$ws.Range("B12").Value2 = 42222
$ws.Range("A13").Value = "Donnerstag"
$ws.Range("B13").Value2 = 42223
$ws.Range("A14").Value = "Freitag"
$ws.Range("B14").Value2 = 42224
$ws.Range("A15").Value = "Samstag"
$ws.Range("B15").Value2 = 42225
$ws.Range("A16").Value = "Sonntag"
$ws.Range("B16").Value2 = 42226
$ws.Range("A17").Value = "Montag"
$ws.Range("B17").Value2 = 42227
$ws.Range("A18").Value = "Dienstag"
$ws.Range("B18").Value2 = 42228
$ws.Range("A19").Value = "Mittwoch"
$ws.Range("B19").Value2 = 42229
$ws.Range("A20").Value = "Donnerstag"
$ws.Range("B20").Value2 = 42230
$ws.Range("A21").Value = "Freitag"
$ws.Range("B21").Value2 = 42231
$ws.Range("A22").Value = "Samstag"
$ws.Range("B22").Value2 = 42232
$ws.Range("A23").Value = "Sonntag"
$ws.Range("B23").Value2 = 42233
$ws.Range("A24").Value = "Montag"
$ws.Range("B24").Value2 = 42234
$ws.Range("A25").Value = "Dienstag"
$ws.Range("B25").Value2 = 42235
$ws.Range("A26").Value = "Mittwoch"
$ws.Range("B26").Value2 = 42236
$ws.Range("A27").Value = "Donnerstag"
$ws.Range("B27").Value2 = 42237
$ws.Range("A28").Value = "Freitag"
$ws.Range("B28").Value2 = 42238
$ws.Range("A29").Value = "Samstag"
$ws.Range("B29").Value2 = 42239
$ws.Range("A30").Value = "Sonntag"
$ws.Range("B30").Value2 = 42240
$ws.Range("A31").Value = "Montag"
$ws.Range("B31").Value2 = 42241
$ws.Range("A32").Value = "Dienstag"
$ws.Range("B32").Value2 = 42242
$ws.Range("A33").Value = "Mittwoch"
$ws.Range("B33").Value2 = 42243
$ws.Range("A34").Value = "Donnerstag"
$ws.Range("B34").Value2 = 42244
$ws.Range("A35").Value = "Freitag"
$ws.Range("B35").Value2 = 42245
$ws.Range("A36").Value = "Samstag"
$ws.Range("B36").Value2 = 42246
